# Cambio cabeceras a castellano
# Translate the table header row from Catalan to Spanish.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Código"
$ws.Range("B1").Value = "Centro"
$ws.Range("C1").Value = "Rég."
$ws.Range("D1").Value = "Dirección"
$ws.Range("E1").Value = "Localidad"
$ws.Range("F1").Value = "Teléfono"

# Move the selection from the old cell (C87) to the header row.
$ws.Range("A1:F1").Select()
